$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20
$ws.Cells.Item(20,10).Value = ""
$ws.Cells.Item(20,1).Value = 111926622
$ws.Cells.Item(20,2).Value = 90658
$ws.Cells.Item(20,3).Value = "'Ovaliderad"
$ws.Cells.Item(20,4).Value = "'NT"
$ws.Cells.Item(20,5).Value = 4361
$ws.Cells.Item(20,6).Value = "'Orange taggsvamp"
$ws.Cells.Item(20,7).Value = "'Hydnellum aurantiacum"
$ws.Cells.Item(20,8).Value = "'(Batsch:Fr.) P.Karst."
$ws.Cells.Item(20,9).Value = ""
$ws.Cells.Item(20,11).Value = ""
$ws.Cells.Item(20,16).Value = "'Upplands Väsby (Upplands Väsby), Upl"
$ws.Cells.Item(20,17).Value = 663452.3464515609
$ws.Cells.Item(20,18).Value = 6602675.90838708
$ws.Cells.Item(20,19).Value = 10
$ws.Cells.Item(20,20).Value = "'Stockholm"
$ws.Cells.Item(20,21).Value = "'Upplands Väsby"
$ws.Cells.Item(20,22).Value = "'Uppland"
$ws.Cells.Item(20,23).Value = "'Ed"
$ws.Cells.Item(20,25).Value = "'2023-09-06"
$ws.Cells.Item(20,26).Value = "'00:00"
$ws.Cells.Item(20,27).Value = "'2023-09-06"
$ws.Cells.Item(20,28).Value = "'00:00"
$ws.Cells.Item(20,30).Value = $false
$ws.Cells.Item(20,31).Value = $false
$ws.Cells.Item(20,33).Value = $false
$ws.Cells.Item(20,46).Value = ""
$ws.Cells.Item(20,49).Value = "'Hans Bärring"
$ws.Cells.Item(20,50).Value = "'Hans Bärring"
$ws.Cells.Item(20,51).Value = ""

# Row 21
$ws.Cells.Item(21,1).Value = 111929648
$ws.Cells.Item(21,2).Value = 90655
$ws.Cells.Item(21,3).Value = "'Ovaliderad"
$ws.Cells.Item(21,4).Value = "'VU"
$ws.Cells.Item(21,5).Value = 150
$ws.Cells.Item(21,6).Value = "'Grangråticka"
$ws.Cells.Item(21,7).Value = "'Boletopsis leucomelaena"
$ws.Cells.Item(21,8).Value = "'(Pers.) Fayod"
$ws.Cells.Item(21,9).Value = "'10"
$ws.Cells.Item(21,10).Value = "'fruktkroppar"
$ws.Cells.Item(21,11).Value = ""
$ws.Cells.Item(21,16).Value = "'Upplands Väsby (Upplands Väsby), Upl"
$ws.Cells.Item(21,17).Value = 663509.4011657666
$ws.Cells.Item(21,18).Value = 6602732.737488487
$ws.Cells.Item(21,19).Value = 10
$ws.Cells.Item(21,20).Value = "'Stockholm"
$ws.Cells.Item(21,21).Value = "'Upplands Väsby"
$ws.Cells.Item(21,22).Value = "'Uppland"
$ws.Cells.Item(21,23).Value = "'Ed"
$ws.Cells.Item(21,25).Value = "'2023-09-06"
$ws.Cells.Item(21,26).Value = "'00:00"
$ws.Cells.Item(21,27).Value = "'2023-09-06"
$ws.Cells.Item(21,28).Value = "'00:00"
$ws.Cells.Item(21,30).Value = $false
$ws.Cells.Item(21,31).Value = $false
$ws.Cells.Item(21,33).Value = $false
$ws.Cells.Item(21,46).Value = ""
$ws.Cells.Item(21,49).Value = "'Hans Bärring"
$ws.Cells.Item(21,50).Value = "'Hans Bärring"
$ws.Cells.Item(21,51).Value = ""

# Row 22
$ws.Cells.Item(22,1).Value = 111927215
$ws.Cells.Item(22,2).Value = 90658
$ws.Cells.Item(22,3).Value = "'Ovaliderad"
$ws.Cells.Item(22,4).Value = "'NT"
$ws.Cells.Item(22,5).Value = 4361
$ws.Cells.Item(22,6).Value = "'Orange taggsvamp"
$ws.Cells.Item(22,7).Value = "'Hydnellum aurantiacum"
$ws.Cells.Item(22,8).Value = "'(Batsch:Fr.) P.Karst."
$ws.Cells.Item(22,9).Value = "'10"
$ws.Cells.Item(22,10).Value = "'fruktkroppar"
$ws.Cells.Item(22,11).Value = ""
$ws.Cells.Item(22,16).Value = "'Upplands Väsby (Upplands Väsby), Upl"
$ws.Cells.Item(22,17).Value = 663485.6413922446
$ws.Cells.Item(22,18).Value = 6602647.390513759
$ws.Cells.Item(22,19).Value = 10
$ws.Cells.Item(22,20).Value = "'Stockholm"
$ws.Cells.Item(22,21).Value = "'Upplands Väsby"
$ws.Cells.Item(22,22).Value = "'Uppland"
$ws.Cells.Item(22,23).Value = "'Ed"
$ws.Cells.Item(22,25).Value = "'2023-09-06"
$ws.Cells.Item(22,26).Value = "'00:00"
$ws.Cells.Item(22,27).Value = "'2023-09-06"
$ws.Cells.Item(22,28).Value = "'00:00"
$ws.Cells.Item(22,30).Value = $false
$ws.Cells.Item(22,31).Value = $false
$ws.Cells.Item(22,33).Value = $false
$ws.Cells.Item(22,46).Value = ""
$ws.Cells.Item(22,49).Value = "'Hans Bärring"
$ws.Cells.Item(22,50).Value = "'Hans Bärring"
$ws.Cells.Item(22,51).Value = ""

# Row 23
$ws.Cells.Item(23,1).Value = 112083737
$ws.Cells.Item(23,2).Value = 98535
$ws.Cells.Item(23,3).Value = "'Ovaliderad"
$ws.Cells.Item(23,4).Value = "'LC"
$ws.Cells.Item(23,5).Value = 222498
$ws.Cells.Item(23,6).Value = "'Blåsippa"
$ws.Cells.Item(23,7).Value = "'Hepatica nobilis"
$ws.Cells.Item(23,8).Value = "'Schreb."
$ws.Cells.Item(23,9).Value = "'200"
$ws.Cells.Item(23,10).Value = "'stjälkar/strån/skott"
$ws.Cells.Item(23,11).Value = ""
$ws.Cells.Item(23,16).Value = "'Upplands Väsby (Upplands Väsby), Upl"
$ws.Cells.Item(23,17).Value = 663545.1917381487
$ws.Cells.Item(23,18).Value = 6602752.072187248
$ws.Cells.Item(23,19).Value = 30
$ws.Cells.Item(23,20).Value = "'Stockholm"
$ws.Cells.Item(23,21).Value = "'Upplands Väsby"
$ws.Cells.Item(23,22).Value = "'Uppland"
$ws.Cells.Item(23,23).Value = "'Ed"
$ws.Cells.Item(23,25).Value = "'2023-09-14"
$ws.Cells.Item(23,26).Value = "'00:00"
$ws.Cells.Item(23,27).Value = "'2023-09-14"
$ws.Cells.Item(23,28).Value = "'00:00"
$ws.Cells.Item(23,30).Value = $false
$ws.Cells.Item(23,31).Value = $false
$ws.Cells.Item(23,33).Value = $false
$ws.Cells.Item(23,46).Value = ""
$ws.Cells.Item(23,49).Value = "'Hans Bärring"
$ws.Cells.Item(23,50).Value = "'Hans Bärring"
$ws.Cells.Item(23,51).Value = ""

# Row 24
$ws.Cells.Item(24,1).Value = 112083804
$ws.Cells.Item(24,2).Value = 98535
$ws.Cells.Item(24,3).Value = "'Ovaliderad"
$ws.Cells.Item(24,4).Value = "'LC"
$ws.Cells.Item(24,5).Value = 222498
$ws.Cells.Item(24,6).Value = "'Blåsippa"
$ws.Cells.Item(24,7).Value = "'Hepatica nobilis"
$ws.Cells.Item(24,8).Value = "'Schreb."
$ws.Cells.Item(24,9).Value = "'300"
$ws.Cells.Item(24,10).Value = "'stjälkar/strån/skott"
$ws.Cells.Item(24,11).Value = ""
$ws.Cells.Item(24,16).Value = "'Upplands Väsby (Upplands Väsby), Upl"
$ws.Cells.Item(24,17).Value = 663571.7306570449
$ws.Cells.Item(24,18).Value = 6602738.498618284
$ws.Cells.Item(24,19).Value = 10
$ws.Cells.Item(24,20).Value = "'Stockholm"
$ws.Cells.Item(24,21).Value = "'Upplands Väsby"
$ws.Cells.Item(24,22).Value = "'Uppland"
$ws.Cells.Item(24,23).Value = "'Ed"
$ws.Cells.Item(24,25).Value = "'2023-09-14"
$ws.Cells.Item(24,26).Value = "'00:00"
$ws.Cells.Item(24,27).Value = "'2023-09-14"
$ws.Cells.Item(24,28).Value = "'00:00"
$ws.Cells.Item(24,30).Value = $false
$ws.Cells.Item(24,31).Value = $false
$ws.Cells.Item(24,33).Value = $false
$ws.Cells.Item(24,46).Value = ""
$ws.Cells.Item(24,49).Value = "'Hans Bärring"
$ws.Cells.Item(24,50).Value = "'Hans Bärring"
$ws.Cells.Item(24,51).Value = ""

# Row 25
$ws.Cells.Item(25,1).Value = 112083991
$ws.Cells.Item(25,2).Value = 98535
$ws.Cells.Item(25,3).Value = "'Ovaliderad"
$ws.Cells.Item(25,4).Value = "'LC"
$ws.Cells.Item(25,5).Value = 222498
$ws.Cells.Item(25,6).Value = "'Blåsippa"
$ws.Cells.Item(25,7).Value = "'Hepatica nobilis"
$ws.Cells.Item(25,8).Value = "'Schreb."
$ws.Cells.Item(25,9).Value = "'300"
$ws.Cells.Item(25,10).Value = "'stjälkar/strån/skott"
$ws.Cells.Item(25,11).Value = ""
$ws.Cells.Item(25,16).Value = "'Upplands Väsby (Upplands Väsby), Upl"
$ws.Cells.Item(25,17).Value = 663568.3519142884
$ws.Cells.Item(25,18).Value = 6602664.1969273
$ws.Cells.Item(25,19).Value = 10
$ws.Cells.Item(25,20).Value = "'Stockholm"
$ws.Cells.Item(25,21).Value = "'Upplands Väsby"
$ws.Cells.Item(25,22).Value = "'Uppland"
$ws.Cells.Item(25,23).Value = "'Ed"
$ws.Cells.Item(25,25).Value = "'2023-09-14"
$ws.Cells.Item(25,26).Value = "'00:00"
$ws.Cells.Item(25,27).Value = "'2023-09-14"
$ws.Cells.Item(25,28).Value = "'00:00"
$ws.Cells.Item(25,30).Value = $false
$ws.Cells.Item(25,31).Value = $false
$ws.Cells.Item(25,33).Value = $false
$ws.Cells.Item(25,46).Value = ""
$ws.Cells.Item(25,49).Value = "'Hans Bärring"
$ws.Cells.Item(25,50).Value = "'Hans Bärring"
$ws.Cells.Item(25,51).Value = ""

# Row 26
$ws.Cells.Item(26,1).Value = 112083905
$ws.Cells.Item(26,2).Value = 98535
$ws.Cells.Item(26,3).Value = "'Ovaliderad"
$ws.Cells.Item(26,4).Value = "'LC"
$ws.Cells.Item(26,5).Value = 222498
$ws.Cells.Item(26,6).Value = "'Blåsippa"
$ws.Cells.Item(26,7).Value = "'Hepatica nobilis"
$ws.Cells.Item(26,8).Value = "'Schreb."
$ws.Cells.Item(26,9).Value = "'400"
$ws.Cells.Item(26,10).Value = "'stjälkar/strån/skott"
$ws.Cells.Item(26,11).Value = "'fullt utvecklade blad"
$ws.Cells.Item(26,16).Value = "'Upplands Väsby (Upplands Väsby), Upl"
$ws.Cells.Item(26,17).Value = 663567.9108240836
$ws.Cells.Item(26,18).Value = 6602721.063539478
$ws.Cells.Item(26,19).Value = 5
$ws.Cells.Item(26,20).Value = "'Stockholm"
$ws.Cells.Item(26,21).Value = "'Upplands Väsby"
$ws.Cells.Item(26,22).Value = "'Uppland"
$ws.Cells.Item(26,23).Value = "'Ed"
$ws.Cells.Item(26,25).Value = "'2023-09-14"
$ws.Cells.Item(26,26).Value = "'00:00"
$ws.Cells.Item(26,27).Value = "'2023-09-14"
$ws.Cells.Item(26,28).Value = "'00:00"
$ws.Cells.Item(26,30).Value = $false
$ws.Cells.Item(26,31).Value = $false
$ws.Cells.Item(26,33).Value = $false
$ws.Cells.Item(26,46).Value = ""
$ws.Cells.Item(26,49).Value = "'Hans Bärring"
$ws.Cells.Item(26,50).Value = "'Hans Bärring"
$ws.Cells.Item(26,51).Value = ""

# Row 27
$ws.Cells.Item(27,1).Value = 112084535
$ws.Cells.Item(27,2).Value = 88899
$ws.Cells.Item(27,3).Value = "'Ovaliderad"
$ws.Cells.Item(27,4).Value = "'NT"
$ws.Cells.Item(27,5).Value = 3286
$ws.Cells.Item(27,6).Value = "'Flattoppad klubbsvamp"
$ws.Cells.Item(27,7).Value = "'Clavariadelphus truncatus"
$ws.Cells.Item(27,8).Value = "'(Quél.) Donk"
$ws.Cells.Item(27,9).Value = "'80"
$ws.Cells.Item(27,10).Value = "'fruktkroppar"
$ws.Cells.Item(27,11).Value = ""
$ws.Cells.Item(27,16).Value = "'Upplands Väsby (Upplands Väsby), Upl"
$ws.Cells.Item(27,17).Value = 663374.2695844367
$ws.Cells.Item(27,18).Value = 6602611.054278261
$ws.Cells.Item(27,19).Value = 10
$ws.Cells.Item(27,20).Value = "'Stockholm"
$ws.Cells.Item(27,21).Value = "'Upplands Väsby"
$ws.Cells.Item(27,22).Value = "'Uppland"
$ws.Cells.Item(27,23).Value = "'Ed"
$ws.Cells.Item(27,25).Value = "'2023-09-14"
$ws.Cells.Item(27,26).Value = "'11:46"
$ws.Cells.Item(27,27).Value = "'2023-09-14"
$ws.Cells.Item(27,28).Value = "'11:46"
$ws.Cells.Item(27,29).Value = "'Sötaktig mild smak (ej bitter)"
$ws.Cells.Item(27,30).Value = $false
$ws.Cells.Item(27,31).Value = $false
$ws.Cells.Item(27,33).Value = $false
$ws.Cells.Item(27,46).Value = ""
$ws.Cells.Item(27,49).Value = "'Hans Bärring"
$ws.Cells.Item(27,50).Value = "'Hans Bärring"
$ws.Cells.Item(27,51).Value = ""

# Row 28
$ws.Cells.Item(28,1).Value = 112084114
$ws.Cells.Item(28,2).Value = 98535
$ws.Cells.Item(28,3).Value = "'Ovaliderad"
$ws.Cells.Item(28,4).Value = "'LC"
$ws.Cells.Item(28,5).Value = 222498
$ws.Cells.Item(28,6).Value = "'Blåsippa"
$ws.Cells.Item(28,7).Value = "'Hepatica nobilis"
$ws.Cells.Item(28,8).Value = "'Schreb."
$ws.Cells.Item(28,9).Value = "'20"
$ws.Cells.Item(28,10).Value = "'plantor/tuvor"
$ws.Cells.Item(28,11).Value = "'fullt utvecklade blad"
$ws.Cells.Item(28,16).Value = "'Upplands Väsby (Upplands Väsby), Upl"
$ws.Cells.Item(28,17).Value = 663576.8087203993
$ws.Cells.Item(28,18).Value = 6602715.356141716
$ws.Cells.Item(28,19).Value = 5
$ws.Cells.Item(28,20).Value = "'Stockholm"
$ws.Cells.Item(28,21).Value = "'Upplands Väsby"
$ws.Cells.Item(28,22).Value = "'Uppland"
$ws.Cells.Item(28,23).Value = "'Ed"
$ws.Cells.Item(28,25).Value = "'2023-09-14"
$ws.Cells.Item(28,26).Value = "'00:00"
$ws.Cells.Item(28,27).Value = "'2023-09-14"
$ws.Cells.Item(28,28).Value = "'00:00"
$ws.Cells.Item(28,30).Value = $false
$ws.Cells.Item(28,31).Value = $false
$ws.Cells.Item(28,33).Value = $false
$ws.Cells.Item(28,46).Value = ""
$ws.Cells.Item(28,49).Value = "'Hans Bärring"
$ws.Cells.Item(28,50).Value = "'Hans Bärring"
$ws.Cells.Item(28,51).Value = ""

# Row 29
$ws.Cells.Item(29,1).Value = 112083958
$ws.Cells.Item(29,2).Value = 98535
$ws.Cells.Item(29,3).Value = "'Ovaliderad"
$ws.Cells.Item(29,4).Value = "'LC"
$ws.Cells.Item(29,5).Value = 222498
$ws.Cells.Item(29,6).Value = "'Blåsippa"
$ws.Cells.Item(29,7).Value = "'Hepatica nobilis"
$ws.Cells.Item(29,8).Value = "'Schreb."
$ws.Cells.Item(29,9).Value = "'10"
$ws.Cells.Item(29,10).Value = "'plantor/tuvor"
$ws.Cells.Item(29,11).Value = "'fullt utvecklade blad"
$ws.Cells.Item(29,16).Value = "'Upplands Väsby (Upplands Väsby), Upl"
$ws.Cells.Item(29,17).Value = 663551.019940288
$ws.Cells.Item(29,18).Value = 6602700.011799707
$ws.Cells.Item(29,19).Value = 5
$ws.Cells.Item(29,20).Value = "'Stockholm"
$ws.Cells.Item(29,21).Value = "'Upplands Väsby"
$ws.Cells.Item(29,22).Value = "'Uppland"
$ws.Cells.Item(29,23).Value = "'Ed"
$ws.Cells.Item(29,25).Value = "'2023-09-14"
$ws.Cells.Item(29,26).Value = "'00:00"
$ws.Cells.Item(29,27).Value = "'2023-09-14"
$ws.Cells.Item(29,28).Value = "'00:00"
$ws.Cells.Item(29,30).Value = $false
$ws.Cells.Item(29,31).Value = $false
$ws.Cells.Item(29,33).Value = $false
$ws.Cells.Item(29,46).Value = ""
$ws.Cells.Item(29,49).Value = "'Hans Bärring"
$ws.Cells.Item(29,50).Value = "'Hans Bärring"
$ws.Cells.Item(29,51).Value = ""

# Row 30
$ws.Cells.Item(30,1).Value = 112084040
$ws.Cells.Item(30,2).Value = 98535
$ws.Cells.Item(30,3).Value = "'Ovaliderad"
$ws.Cells.Item(30,4).Value = "'LC"
$ws.Cells.Item(30,5).Value = 222498
$ws.Cells.Item(30,6).Value = "'Blåsippa"
$ws.Cells.Item(30,7).Value = "'Hepatica nobilis"
$ws.Cells.Item(30,8).Value = "'Schreb."
$ws.Cells.Item(30,9).Value = ""
$ws.Cells.Item(30,11).Value = ""
$ws.Cells.Item(30,16).Value = "'Upplands Väsby (Upplands Väsby), Upl"
$ws.Cells.Item(30,17).Value = 663584.9559231531
$ws.Cells.Item(30,18).Value = 6602703.52117154
$ws.Cells.Item(30,19).Value = 10
$ws.Cells.Item(30,20).Value = "'Stockholm"
$ws.Cells.Item(30,21).Value = "'Upplands Väsby"
$ws.Cells.Item(30,22).Value = "'Uppland"
$ws.Cells.Item(30,23).Value = "'Ed"
$ws.Cells.Item(30,25).Value = "'2023-09-14"
$ws.Cells.Item(30,26).Value = "'00:00"
$ws.Cells.Item(30,27).Value = "'2023-09-14"
$ws.Cells.Item(30,28).Value = "'00:00"
$ws.Cells.Item(30,30).Value = $false
$ws.Cells.Item(30,31).Value = $false
$ws.Cells.Item(30,33).Value = $false
$ws.Cells.Item(30,46).Value = ""
$ws.Cells.Item(30,49).Value = "'Hans Bärring"
$ws.Cells.Item(30,50).Value = "'Hans Bärring"
$ws.Cells.Item(30,51).Value = ""

